# Fruta / hortaliza, semanal
# Insert two new weekly rows at the top of the Repollo / Macroferia Regional de
# Talca block (row 293), pushing the existing data down by two rows, then
# populate the two newly-inserted rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 293:319 down to 295:321.
$ws.Rows("293:294").Insert()

# New row 293: Primera
$ws.Range("A293").Value = 5
$ws.Range("B293").Value = "Macroferia Regional de Talca"
$ws.Range("C293").Value = "Maule"
$ws.Range("D293").Value = 44769
$ws.Range("E293").Value = 7
$ws.Range("F293").Value = 100112006
$ws.Range("G293").Value = "Repollo"
$ws.Range("H293").Value = "Crespo record"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 3000
$ws.Range("K293").Value = 1000
$ws.Range("L293").Value = 1000
$ws.Range("M293").Value = 1000
$ws.Range("N293").Value = "$/unidad"
$ws.Range("O293").Value = "Región del Maule"
$ws.Range("P293").Value = 1000
$ws.Range("Q293").Value = 1
$ws.Range("R293").Value = "Hortaliza"

# New row 294: Segunda
$ws.Range("A294").Value = 5
$ws.Range("B294").Value = "Macroferia Regional de Talca"
$ws.Range("C294").Value = "Maule"
$ws.Range("D294").Value = 44769
$ws.Range("E294").Value = 7
$ws.Range("F294").Value = 100112006
$ws.Range("G294").Value = "Repollo"
$ws.Range("H294").Value = "Crespo record"
$ws.Range("I294").Value = "Segunda"
$ws.Range("J294").Value = 3000
$ws.Range("K294").Value = 800
$ws.Range("L294").Value = 800
$ws.Range("M294").Value = 800
$ws.Range("N294").Value = "$/unidad"
$ws.Range("O294").Value = "Región del Maule"
$ws.Range("P294").Value = 800
$ws.Range("Q294").Value = 1
$ws.Range("R294").Value = "Hortaliza"
